# Update cryptos price/volume table to reflect latest scrape.
# Generated from the diff of cryptos.xlsx (rows 2-51, columns B-E).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.531.66"
$ws.Range("E2").Value = "  -3.07%  "
$ws.Range("D3").Value = "3.003.63"
$ws.Range("E3").Value = "  -2.64%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'544.07"
$ws.Range("E5").Value = "  -0.63%  "
$ws.Range("D6").Value = "'130.77"
$ws.Range("E6").Value = "  -6.33%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "2.997.65"
$ws.Range("E8").Value = "  -2.62%  "
$ws.Range("D9").Value = "'0.489"
$ws.Range("E9").Value = "  -1.74%  "
$ws.Range("B10").Value = "Toncoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D10").Value = "'6.00"
$ws.Range("E10").Value = "  -6.31%  "
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").Value = "'0.145"
$ws.Range("E11").Value = "  -8.66%  "
$ws.Range("D12").Value = "'0.443"
$ws.Range("E12").Value = "  -3.37%  "
$ws.Range("E13").Value = "  -3.14%  "
$ws.Range("D14").Value = "'33.86"
$ws.Range("E14").Value = "  -3.38%  "
$ws.Range("D15").Value = "3.485.12"
$ws.Range("E15").Value = "  -2.82%  "
$ws.Range("D16").Value = "61.653.19"
$ws.Range("E16").Value = "  -2.94%  "
$ws.Range("D17").Value = "'0.110"
$ws.Range("E17").Value = "  -2.89%  "
$ws.Range("D18").Value = "2.998.18"
$ws.Range("E18").Value = "  -2.77%  "
$ws.Range("E19").Value = "  -1.15%  "
$ws.Range("D20").Value = "'479.65"
$ws.Range("E20").Value = "  +0.77%  "
$ws.Range("E21").Value = "  -2.75%  "
$ws.Range("D22").Value = "'0.662"
$ws.Range("E22").Value = "  -5.97%  "
$ws.Range("D23").Value = "'6.95"
$ws.Range("E23").Value = "  -2.06%  "
$ws.Range("D24").Value = "'80.79"
$ws.Range("E24").Value = "  +2.55%  "
$ws.Range("D25").Value = "'11.98"
$ws.Range("E25").Value = "  -2.23%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").Value = "'2.69"
$ws.Range("E27").Value = "  -1.63%  "
$ws.Range("D28").Value = "'7.61"
$ws.Range("E28").Value = "  -4.43%  "
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("E30").Value = "  -0.14%  "
$ws.Range("D31").Value = "'25.51"
$ws.Range("E31").Value = "  -3.05%  "
$ws.Range("E32").Value = "  -3.67%  "
$ws.Range("D33").Value = "'2.32"
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("E34").Value = "  -0.38%  "
$ws.Range("D35").Value = "'54.69"
$ws.Range("E35").Value = "  -7.22%  "
$ws.Range("E36").Value = "  -3.28%  "
$ws.Range("D37").Value = "'444.15"
$ws.Range("E37").Value = "  -9.37%  "
$ws.Range("D38").Value = "3.128.04"
$ws.Range("E38").Value = "  -4.17%  "
$ws.Range("E39").Value = "  -1.14%  "
$ws.Range("E40").Value = "  -5.78%  "
$ws.Range("D41").Value = "'0.116"
$ws.Range("E41").Value = "  -1.55%  "
$ws.Range("D42").Value = "'8.05"
$ws.Range("E42").Value = "  -1.43%  "
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("D44").Value = "'2.36"
$ws.Range("E44").Value = "  -9.65%  "
$ws.Range("D45").Value = "'25.60"
$ws.Range("E45").Value = "  +0.41%  "
$ws.Range("E46").Value = "  -4.97%  "
$ws.Range("E47").Value = "  -1.71%  "
$ws.Range("D48").Value = "'1.93"
$ws.Range("E48").Value = "  -4.73%  "
$ws.Range("B49").Value = "BitgetToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/q7gMmMdLb+bitgettoken-bgb"
$ws.Range("D49").Value = "'1.30"
$ws.Range("E49").Value = "  +9.15%  "
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").Value = "'113.98"
$ws.Range("E50").Value = "  -8.51%  "
$ws.Range("D51").Value = "0.0₃0482"
$ws.Range("E51").Value = "  -9.34%  "
